$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    if ($value -eq "") {
        # A bare "" clears the cell entirely (becomes a true empty/null
        # cell). The source file instead keeps an explicit empty *string*
        # cell there, so force that with a leading-quote empty literal.
        $cell.Value = "'"
    } else {
        $cell.Value = $value
    }
}

# Structural change: a new row is inserted at row 7 (the old "Сумма:" totals
# row ends up re-created there with updated numbers; row 6 becomes blank).
$ws.Rows.Item(7).Insert()

# Row 2
Set-TextValue 2 2 "2.1 - 2.1"

# Row 3
Set-TextValue 3 2 "1.2 - 2.3"
Set-TextValue 3 3 "30"
Set-TextValue 3 8 "9808"
Set-TextValue 3 9 "jkhkjh"
Set-TextValue 3 11 "87"

# Row 4
Set-TextValue 4 2 "1.1 - 1.1"
Set-TextValue 4 6 "98098080"
Set-TextValue 4 7 ""
Set-TextValue 4 8 ""
Set-TextValue 4 10 ""
Set-TextValue 4 11 ""

# Row 5
Set-TextValue 5 2 "27.7 - 27.7"
Set-TextValue 5 3 "1"
Set-TextValue 5 6 "asdfsa23"
Set-TextValue 5 7 "234"
Set-TextValue 5 8 "234"
Set-TextValue 5 10 "2df"
Set-TextValue 5 11 "234"

# Row 6 (old totals row content is cleared - the totals move to the new row 7)
Set-TextValue 6 2 ""
Set-TextValue 6 6 ""
Set-TextValue 6 7 ""
Set-TextValue 6 8 ""
Set-TextValue 6 11 ""

# Row 7 (new totals row) - column A keeps the same bold/bordered/centered
# look used by the other numbered rows (A2:A6).
$aCell = $ws.Cells.Item(7, 1)
$aCell.Font.Bold = $true
$aCell.HorizontalAlignment = -4108   # xlCenter
$aCell.VerticalAlignment = -4160     # xlTop
$aCell.Borders.LineStyle = 1
$aCell.Borders.Weight = 2
$aCell.Value = 5
Set-TextValue 7 2 "Сумма:"
Set-TextValue 7 3 ""
Set-TextValue 7 4 ""
Set-TextValue 7 5 ""
Set-TextValue 7 6 "98098080"
Set-TextValue 7 7 "234"
Set-TextValue 7 8 "10042"
Set-TextValue 7 9 ""
Set-TextValue 7 10 ""
Set-TextValue 7 11 "321"

$wb.Save()
